$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect it so the cell values below can be edited,
# then re-protect it again once the edits are made.
$ws.Unprotect()

# Update the "as of" date in the confidentiality footnote: 2021-05-06 -> 2021-05-07
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) figures for each holding (rows 2-12)
$ws.Range("D2").Value = 0.03054491706223923
$ws.Range("E2").Value = 0.01272352132049503

$ws.Range("D3").Value = 0.02363545316794995
$ws.Range("E3").Value = 0.003285613705702906

$ws.Range("D4").Value = 0.05182363465063392
$ws.Range("E4").Value = 0.001393728222996771

$ws.Range("D5").Value = 0.1370463750187209
$ws.Range("E5").Value = 0.00732421875

$ws.Range("D6").Value = 0.03138309303930376
$ws.Range("E6").Value = 0.02039381153305198

$ws.Range("D7").Value = 0.119728162842651
$ws.Range("E7").Value = 0.009889834752128124

$ws.Range("D8").Value = 0.1026523631499508
$ws.Range("E8").Value = 0.01089918256130784

$ws.Range("D9").Value = 0.02939682132795541
$ws.Range("E9").Value = 0.008804258804258724

$ws.Range("D10").Value = 0.1267295640906623
$ws.Range("E10").Value = 0.005587632706276713

$ws.Range("D11").Value = 0.2443739709418953
$ws.Range("E11").Value = 0.01000183519911912

$ws.Range("D12").Value = 0.1026856447080375
$ws.Range("E12").Value = 0.01023363583703429

# Row 13 "Total": only the Percent Change total changes; Weight stays 1
$ws.Range("E13").Value = 0.008947192198991649

# Restore worksheet protection
$ws.Protect()
